$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4800.8
$ws.Range("I40").Value = 4750.5
$ws.Range("K40").Value = 4750.5
$ws.Range("M40").Value = -4575.5
$ws.Range("H80").Value = 727.625
$ws.Range("I80").Value = 210.13333
$ws.Range("J80").Value = 1590.1111
$ws.Range("K80").Value = 630.39999
$ws.Range("L80").Value = 4770.3333
$ws.Range("M80").Value = 367.60001
$ws.Range("N80").Value = -6766.3333
$ws.Range("H83").Value = 727.625
$ws.Range("I83").Value = 210.13333
$ws.Range("J83").Value = 1590.1111
$ws.Range("K83").Value = 1891.19997
$ws.Range("L83").Value = 14310.9999
$ws.Range("M83").Value = 3100.80003
$ws.Range("N83").Value = -24294.9999
$ws.Range("H103").Value = 717.1
$ws.Range("I103").Value = 749
$ws.Range("J103").Value = 685.2
$ws.Range("K103").Value = 2247
$ws.Range("L103").Value = 2055.6
$ws.Range("M103").Value = -1661
$ws.Range("N103").Value = -3227.6
$ws.Range("H116").Value = 5338.5713
$ws.Range("I116").Value = 5293.643
$ws.Range("J116").Value = 5428.4287
$ws.Range("K116").Value = 5293.643
$ws.Range("L116").Value = 5428.4287
$ws.Range("M116").Value = -1851.643
$ws.Range("N116").Value = -12312.4287
$ws.Range("H137").Value = 5117.694
$ws.Range("I137").Value = 2395.6
$ws.Range("K137").Value = 7186.799999999999
$ws.Range("M137").Value = -4636.799999999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23176.926
$ws.Range("I32").Value = 23589.963
$ws.Range("K32").Value = 23589.963
$ws.Range("M32").Value = -23302.963
$ws.Range("H61").Value = 6806084
$ws.Range("I61").Value = 7578114.5
$ws.Range("J61").Value = 12214
$ws.Range("K61").Value = 7578114.5
$ws.Range("L61").Value = 12214
$ws.Range("M61").Value = -7577902.5
$ws.Range("N61").Value = -12638
$ws.Range("H97").Value = 3794849
$ws.Range("I97").Value = 5291497
$ws.Range("J97").Value = 302670
$ws.Range("K97").Value = 5291497
$ws.Range("L97").Value = 302670
$ws.Range("M97").Value = -5291001
$ws.Range("N97").Value = -303662
$ws.Range("H122").Value = 3129.1428
$ws.Range("I122").Value = 2985.2307
$ws.Range("K122").Value = 8955.6921
$ws.Range("M122").Value = -6505.6921
$ws.Range("H132").Value = 3640881
$ws.Range("I132").Value = 4351033.5
$ws.Range("J132").Value = 11212.556
$ws.Range("K132").Value = 13053100.5
$ws.Range("L132").Value = 33637.66800000001
$ws.Range("M132").Value = -13050570.5
$ws.Range("N132").Value = -38697.66800000001
$ws.Range("H136").Value = 6806084
$ws.Range("I136").Value = 7578114.5
$ws.Range("J136").Value = 12214
$ws.Range("K136").Value = 22734343.5
$ws.Range("L136").Value = 36642
$ws.Range("M136").Value = -22731793.5
$ws.Range("N136").Value = -41742

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1927.3077
$ws.Range("J94").Value = 1968.2222
$ws.Range("L94").Value = 1968.2222
$ws.Range("N94").Value = -2870.2222
$ws.Range("H134").Value = 5869.1724
$ws.Range("I134").Value = 5081.7036
$ws.Range("J134").Value = 16500
$ws.Range("K134").Value = 15245.1108
$ws.Range("L134").Value = 49500
$ws.Range("M134").Value = -12710.1108
$ws.Range("N134").Value = -54570

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 456.375
$ws.Range("I33").Value = 387.5
$ws.Range("J33").Value = 525.25
$ws.Range("K33").Value = 2325
$ws.Range("L33").Value = 3151.5
$ws.Range("M33").Value = -2042
$ws.Range("N33").Value = -3717.5
$ws.Range("H34").Value = 1145.6296
$ws.Range("J34").Value = 2306.0833
$ws.Range("L34").Value = 6918.249899999999
$ws.Range("N34").Value = -7086.249899999999
$ws.Range("H138").Value = 402539.62
$ws.Range("I138").Value = 1887.8125
$ws.Range("J138").Value = 1114809.5
$ws.Range("K138").Value = 5663.4375
$ws.Range("L138").Value = 3344428.5
$ws.Range("M138").Value = -523.4375
$ws.Range("N138").Value = -3354708.5
$ws.Range("H140").Value = 22728704
$ws.Range("I140").Value = 45455404
$ws.Range("J140").Value = 2006
$ws.Range("K140").Value = 136366212
$ws.Range("L140").Value = 6018
$ws.Range("M140").Value = -136361032
$ws.Range("N140").Value = -16378

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5734.2085
$ws.Range("I70").Value = 5089.4116
$ws.Range("J70").Value = 7300.143
$ws.Range("K70").Value = 5089.4116
$ws.Range("L70").Value = 7300.143
$ws.Range("M70").Value = -4819.4116
$ws.Range("N70").Value = -7840.143
$ws.Range("H73").Value = 5734.2085
$ws.Range("I73").Value = 5089.4116
$ws.Range("J73").Value = 7300.143
$ws.Range("K73").Value = 5089.4116
$ws.Range("L73").Value = 7300.143
$ws.Range("M73").Value = -4153.4116
$ws.Range("N73").Value = -9172.143
$ws.Range("H80").Value = 6491.231
$ws.Range("I80").Value = 5605.375
$ws.Range("J80").Value = 7908.6
$ws.Range("K80").Value = 5605.375
$ws.Range("L80").Value = 7908.6
$ws.Range("M80").Value = -4607.375
$ws.Range("N80").Value = -9904.6
$ws.Range("H83").Value = 6491.231
$ws.Range("I83").Value = 5605.375
$ws.Range("J83").Value = 7908.6
$ws.Range("K83").Value = 28026.875
$ws.Range("L83").Value = 39543
$ws.Range("M83").Value = -23034.875
$ws.Range("N83").Value = -49527
$ws.Range("H87").Value = 36104.145
$ws.Range("J87").Value = 36104.145
$ws.Range("L87").Value = 36104.145
$ws.Range("N87").Value = -38600.145
$ws.Range("H90").Value = 36104.145
$ws.Range("J90").Value = 36104.145
$ws.Range("L90").Value = 108312.435
$ws.Range("N90").Value = -120792.435
$ws.Range("H94").Value = 25000
$ws.Range("J94").Value = 25000
$ws.Range("L94").Value = 25000
$ws.Range("N94").Value = -26352
$ws.Range("H107").Value = 720.1667
$ws.Range("I107").Value = 1023.1667
$ws.Range("J107").Value = 417.16666
$ws.Range("K107").Value = 1023.1667
$ws.Range("L107").Value = 417.16666
$ws.Range("M107").Value = 896.8333
$ws.Range("N107").Value = -4257.16666
$ws.Range("H132").Value = 3985.8718
$ws.Range("I132").Value = 2597.1765
$ws.Range("K132").Value = 7791.529500000001
$ws.Range("M132").Value = -5261.529500000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 729.24
$ws.Range("I16").Value = 624.2083
$ws.Range("K16").Value = 624.2083
$ws.Range("M16").Value = -454.2083
$ws.Range("H55").Value = 339.96155
$ws.Range("J55").Value = 346
$ws.Range("L55").Value = 346
$ws.Range("N55").Value = -692
$ws.Range("H82").Value = 1994.6154
$ws.Range("I82").Value = 1879.5
$ws.Range("K82").Value = 1879.5
$ws.Range("M82").Value = -1518.5
$ws.Range("H85").Value = 1994.6154
$ws.Range("I85").Value = 1879.5
$ws.Range("K85").Value = 1879.5
$ws.Range("M85").Value = -631.5
$ws.Range("H122").Value = 166669580
$ws.Range("I122").Value = 166669580
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 500008740
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -500006290
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 2303949.8
$ws.Range("I136").Value = 4766056.5
$ws.Range("J136").Value = 5983.467
$ws.Range("K136").Value = 14298169.5
$ws.Range("L136").Value = 17950.401
$ws.Range("M136").Value = -14295619.5
$ws.Range("N136").Value = -23050.401

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1475
$ws.Range("I107").Value = 1174.8125
$ws.Range("K107").Value = 3524.4375
$ws.Range("M107").Value = -1604.4375
$ws.Range("H122").Value = 3850.4348
$ws.Range("I122").Value = 3845.476
$ws.Range("K122").Value = 11536.428
$ws.Range("M122").Value = -9086.428
$ws.Range("H132").Value = 5833.4287
$ws.Range("I132").Value = 2575.2
$ws.Range("K132").Value = 7725.599999999999
$ws.Range("M132").Value = -5195.599999999999
$ws.Range("H136").Value = 4610803
$ws.Range("I136").Value = 5716028
$ws.Range("K136").Value = 17148084
$ws.Range("M136").Value = -17145534
